$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 5

$ws.Cells.Item($row, 1).Value = 42607.890335648146
$ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item($row, 2).Value = 8
$ws.Cells.Item($row, 3).Value = 57
$ws.Cells.Item($row, 4).Value = 42
$ws.Cells.Item($row, 5).Value = 54
$ws.Cells.Item($row, 6).Value = 45
$ws.Cells.Item($row, 7).Value = 19770
$ws.Cells.Item($row, 8).Value = 15626
$ws.Cells.Item($row, 9).Value = 2718
$ws.Cells.Item($row, 10).Value = 288
$ws.Cells.Item($row, 11).Value = 210
$ws.Cells.Item($row, 12).Value = 12
$ws.Cells.Item($row, 13).Value = 10
$ws.Cells.Item($row, 14).Value = "Noun"
